# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the row corresponding to 904b8907-1afe-4039-9b69-7c5f782d8f3c.md now that a
# fresh xliff handoff was generated for it.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-28 12:41:33"

# --- "zh-cn" sheet: column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-28 12:41:29"

# --- "de-de" sheet: column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-28 12:41:33"
